$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated data values for rows 2-9 (columns A:D)
$data = @(
    @(3, 1, 5, 5),
    @(2, 3, 5, 5),
    @(5, 6, 5, 5),
    @(7, 6, 10, 11),
    @(1, 7, 5, 7),
    @(4, 7, 12, 13),
    @(6, 7, 18, 19),
    @(7, 7, 24, 25)
)

$rowIndex = 2
foreach ($rowValues in $data) {
    $ws.Cells.Item($rowIndex, 1).Value = $rowValues[0]
    $ws.Cells.Item($rowIndex, 2).Value = $rowValues[1]
    $ws.Cells.Item($rowIndex, 3).Value = $rowValues[2]
    $ws.Cells.Item($rowIndex, 4).Value = $rowValues[3]
    $rowIndex++
}
